# Fruta / hortaliza, semanal
#
# A new weekly price record for "Ajo" (Macroferia Regional de Talca) is
# inserted as row 426, pushing the existing rows 426-458 down to 427-459.
# The sheet's used range grows from A1:R458 to A1:R459 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 426 (shifts 426:458 -> 427:459, keeps formatting
# of surrounding rows, e.g. the date style on column D).
$ws.Rows.Item(426).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A426").Value = 5
$ws.Range("B426").Value = "Macroferia Regional de Talca"
$ws.Range("C426").Value = "Maule"
$ws.Range("D426").Value = 45013
$ws.Range("E426").Value = 7
$ws.Range("F426").Value = 100112003
$ws.Range("G426").Value = "Ajo"
$ws.Range("H426").Value = "Chino"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 200
$ws.Range("K426").Value = 19000
$ws.Range("L426").Value = 19000
$ws.Range("M426").Value = 19000
$ws.Range("N426").Value = "$/caja 10 kilos"
$ws.Range("O426").Value = "China"
$ws.Range("P426").Value = 1900
$ws.Range("Q426").Value = 10
$ws.Range("R426").Value = "Hortaliza"
